$wb = $excel.ActiveWorkbook

# --- Rename sheets (Monthly -> Weekly) ---
$wsGeneralTax = $wb.Worksheets.Item("GeneralTaxRateMonthly")
$wsGeneralTax.Name = "GeneralTaxRateWeekly"

$wsProcessPayroll = $wb.Worksheets.Item("ProcessPayrollForMonthlyTax")
$wsProcessPayroll.Name = "ProcessPayrollForWeeklyTax"

$wsTestReports = $wb.Worksheets.Item("TestReports")
$wsFirst = $wb.Worksheets.Item("first")

# --- Update "DO NOT TOUCH AUTOMATION EMP 105" -> "... EMP 107" everywhere ---
$wsGeneralTax.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsProcessPayroll.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsTestReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

# --- "first" sheet: update the referenced sheet names in column A ---
$wsFirst.Range("A3").Value = "GeneralTaxRateWeekly"
$wsFirst.Range("A4").Value = "ProcessPayrollForWeeklyTax"

# --- Update selections per sheet (non-active sheets first) ---
$wsGeneralTax.Range("A2").Select()
$wsProcessPayroll.Range("B2").Select()
$wsTestReports.Range("B2").Select()

# --- "first" sheet becomes the active tab with A3 selected (select last) ---
$wsFirst.Activate()
$wsFirst.Range("A3").Select()
